# Update F-column "views/likes" counters on the 展览 (Exhibition) sheet
# and the corresponding rows on the 全部类型 (All types) aggregate sheet,
# matching the regenerated data snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# Mapping: row -> new value, for the 展览 sheet
$sheet1Updates = @{
    2  = 615
    3  = 282
    6  = 415
    10 = 257
    11 = 6998
    15 = 557
    16 = 378
    19 = 18
    25 = 1047
    27 = 32
    28 = 1976
    29 = 546
}

foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

# Mapping: row -> new value, for the 全部类型 sheet
$sheet4Updates = @{
    3  = 615
    4  = 282
    8  = 415
    12 = 257
    13 = 6998
    18 = 557
    19 = 378
    23 = 18
    35 = 1047
    37 = 32
    38 = 1976
    39 = 546
}

foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
